$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1238.4242
$ws.Range("I46").Value = 417
$ws.Range("J46").Value = 1291.4193
$ws.Range("K46").Value = 1251
$ws.Range("L46").Value = 3874.2579
$ws.Range("M46").Value = -1132
$ws.Range("N46").Value = -4112.257900000001

$ws.Range("H60").Value = 1238.4242
$ws.Range("I60").Value = 417
$ws.Range("J60").Value = 1291.4193
$ws.Range("K60").Value = 1251
$ws.Range("L60").Value = 3874.2579
$ws.Range("M60").Value = -767
$ws.Range("N60").Value = -4842.257900000001

$ws.Range("H76").Value = 3515.4546
$ws.Range("I76").Value = 2981.4285
$ws.Range("K76").Value = 2981.4285
$ws.Range("M76").Value = -2666.4285

$ws.Range("H79").Value = 3515.4546
$ws.Range("I79").Value = 2981.4285
$ws.Range("K79").Value = 2981.4285
$ws.Range("M79").Value = -1889.4285

$ws.Range("H103").Value = 705
$ws.Range("I103").Value = 589
$ws.Range("J103").Value = 763
$ws.Range("K103").Value = 1767
$ws.Range("L103").Value = 2289
$ws.Range("M103").Value = -1181
$ws.Range("N103").Value = -3461

$ws.Range("H132").Value = 8936089
$ws.Range("I132").Value = 9267037
$ws.Range("J132").Value = 506
$ws.Range("K132").Value = 27801111
$ws.Range("L132").Value = 1518
$ws.Range("M132").Value = -27798581
$ws.Range("N132").Value = -6578

$ws.Range("H135").Value = 541.29034
$ws.Range("I135").Value = 562.2069
$ws.Range("K135").Value = 5059.8621
$ws.Range("M135").Value = -2524.8621

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2187.2222
$ws.Range("I45").Value = 1889.2778
$ws.Range("J45").Value = 2783.111
$ws.Range("K45").Value = 1889.2778
$ws.Range("L45").Value = 2783.111
$ws.Range("M45").Value = -1512.2778
$ws.Range("N45").Value = -3537.111

$ws.Range("H110").Value = 17277338
$ws.Range("I110").Value = 23858118
$ws.Range("K110").Value = 23858118
$ws.Range("M110").Value = -23856073

$ws.Range("H132").Value = 2466.818
$ws.Range("I132").Value = 3161.3142
$ws.Range("J132").Value = 1251.45
$ws.Range("K132").Value = 9483.942599999998
$ws.Range("L132").Value = 3754.35
$ws.Range("M132").Value = -6953.942599999998
$ws.Range("N132").Value = -8814.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 18000
$ws.Range("J18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("N18").Value = -19058

$ws.Range("H80").Value = 1641
$ws.Range("I80").Value = 755
$ws.Range("J80").Value = 2124.2727
$ws.Range("K80").Value = 755
$ws.Range("L80").Value = 2124.2727
$ws.Range("M80").Value = 243
$ws.Range("N80").Value = -4120.2727

$ws.Range("H83").Value = 1641
$ws.Range("I83").Value = 755
$ws.Range("J83").Value = 2124.2727
$ws.Range("K83").Value = 3775
$ws.Range("L83").Value = 10621.3635
$ws.Range("M83").Value = 1217
$ws.Range("N83").Value = -20605.3635

$ws.Range("H134").Value = 2610.0435
$ws.Range("I134").Value = 2744.1667
$ws.Range("J134").Value = 2127.2
$ws.Range("K134").Value = 8232.500100000001
$ws.Range("L134").Value = 6381.599999999999
$ws.Range("M134").Value = -5697.500100000001
$ws.Range("N134").Value = -11451.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15836.723
$ws.Range("I31").Value = 36794.75
$ws.Range("J31").Value = 2499.7954
$ws.Range("K31").Value = 36794.75
$ws.Range("L31").Value = 2499.7954
$ws.Range("M31").Value = -36499.75
$ws.Range("N31").Value = -3089.7954

$ws.Range("H34").Value = 15836.723
$ws.Range("I34").Value = 36794.75
$ws.Range("J34").Value = 2499.7954
$ws.Range("K34").Value = 36794.75
$ws.Range("L34").Value = 2499.7954
$ws.Range("M34").Value = -36592.75
$ws.Range("N34").Value = -2903.7954

$ws.Range("H105").Value = 973.64
$ws.Range("I105").Value = 949.5714
$ws.Range("J105").Value = 1100
$ws.Range("K105").Value = 949.5714
$ws.Range("L105").Value = 1100
$ws.Range("M105").Value = 797.4286
$ws.Range("N105").Value = -4594

$ws.Range("H107").Value = 628.2143
$ws.Range("I107").Value = 560.95
$ws.Range("J107").Value = 796.375
$ws.Range("K107").Value = 560.95
$ws.Range("L107").Value = 796.375
$ws.Range("M107").Value = 1359.05
$ws.Range("N107").Value = -4636.375

$ws.Range("H132").Value = 2067.2104
$ws.Range("I132").Value = 2161.0938
$ws.Range("J132").Value = 1566.5
$ws.Range("K132").Value = 6483.2814
$ws.Range("L132").Value = 4699.5
$ws.Range("M132").Value = -3953.2814
$ws.Range("N132").Value = -9759.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 450
$ws.Range("J29").Value = 450
$ws.Range("L29").Value = 1350
$ws.Range("N29").Value = -1904

$ws.Range("H46").Value = 125180.5
$ws.Range("I46").Value = 240.66667
$ws.Range("K46").Value = 722.00001
$ws.Range("M46").Value = -631.00001

$ws.Range("H60").Value = 532
$ws.Range("I60").Value = 415
$ws.Range("K60").Value = 1245
$ws.Range("M60").Value = -994

$ws.Range("H68").Value = 1977.9324
$ws.Range("I68").Value = 1275.7273
$ws.Range("J68").Value = 2543.122
$ws.Range("K68").Value = 3827.1819
$ws.Range("L68").Value = 7629.366
$ws.Range("M68").Value = -3016.1819
$ws.Range("N68").Value = -9251.366

$ws.Range("H69").Value = 1651.0834
$ws.Range("I69").Value = 400
$ws.Range("J69").Value = 1764.8182
$ws.Range("K69").Value = 1200
$ws.Range("L69").Value = 5294.4546
$ws.Range("M69").Value = -389
$ws.Range("N69").Value = -6916.4546

$ws.Range("H71").Value = 1977.9324
$ws.Range("I71").Value = 1275.7273
$ws.Range("J71").Value = 2543.122
$ws.Range("K71").Value = 11481.5457
$ws.Range("L71").Value = 22888.098
$ws.Range("M71").Value = -7425.545700000001
$ws.Range("N71").Value = -31000.098

$ws.Range("H72").Value = 1651.0834
$ws.Range("I72").Value = 400
$ws.Range("J72").Value = 1764.8182
$ws.Range("K72").Value = 3600
$ws.Range("L72").Value = 15883.3638
$ws.Range("M72").Value = 456
$ws.Range("N72").Value = -23995.3638

$ws.Range("H113").Value = 559.4865
$ws.Range("J113").Value = 591.6667
$ws.Range("L113").Value = 1775.0001
$ws.Range("N113").Value = -6115.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6497.3335
$ws.Range("I122").Value = 5859.273
$ws.Range("J122").Value = 8252
$ws.Range("K122").Value = 17577.819
$ws.Range("L122").Value = 24756
$ws.Range("M122").Value = -15127.819
$ws.Range("N122").Value = -29656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 377.4
$ws.Range("J22").Value = 429.66666
$ws.Range("L22").Value = 429.66666
$ws.Range("N22").Value = -1019.66666

$ws.Range("H27").Value = 377.4
$ws.Range("J27").Value = 429.66666
$ws.Range("L27").Value = 429.66666
$ws.Range("N27").Value = -643.66666

$ws.Range("H136").Value = 2130
$ws.Range("I136").Value = 1841.875
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 5525.625
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -2975.625
$ws.Range("N136").Value = -16100.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 911.8889
$ws.Range("J122").Value = 870
$ws.Range("L122").Value = 2610
$ws.Range("N122").Value = -7510

$ws.Range("H132").Value = 2320.2327
$ws.Range("I132").Value = 2689.0334
$ws.Range("J132").Value = 1469.1538
$ws.Range("K132").Value = 8067.100199999999
$ws.Range("L132").Value = 4407.4614
$ws.Range("M132").Value = -5537.100199999999
$ws.Range("N132").Value = -9467.4614

$ws.Range("H136").Value = 897.25806
$ws.Range("I136").Value = 529.7917
$ws.Range("J136").Value = 2157.1428
$ws.Range("K136").Value = 1589.3751
$ws.Range("L136").Value = 6471.428400000001
$ws.Range("M136").Value = 960.6249
$ws.Range("N136").Value = -11571.4284
